{"js": "// Assignment # 9 Part 2\n// Replace the single \"College major by Pay.\" paragraph with four new\n// paragraphs of analysis notes, keeping the \"_GoBack\" bookmark anchored\n// right before the final \".\" of the last paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document's only (first) paragraph currently holds \"College major by Pay.\"\nconst targetPara = paragraphs.items[0];\n\n// Insert the three new leading paragraphs directly above it, in document order.\ntargetPara.insertParagraph(\n  \"Anova analysis of median pay by Major Category & Unemployment by major category.\",\n  \"Before\"\n);\ntargetPara.insertParagraph(\"Tukey SD based on Anova.\", \"Before\");\ntargetPara.insertParagraph(\n  \"Bayes prediction of mean pay and unemployment by major category.\",\n  \"Before\"\n);\nawait context.sync();\n\n// Within the original paragraph, swap its trailing text (after the bookmark)\n// down to just \".\", and its leading text (before the bookmark) to the new\n// \"Linear regression ...\" sentence - preserving the bookmark's position.\nconst tailResults = targetPara.search(\"ollege major by Pay.\", { matchCase: true });\ntailResults.load(\"items\");\nawait context.sync();\ntailResults.items[0].insertText(\".\", \"Replace\");\nawait context.sync();\n\nconst headResults = targetPara.search(\"C\", { matchCase: true });\nheadResults.load(\"items\");\nawait context.sync();\nheadResults.items[0].insertText(\n  \"Linear regression of median pay and unemployment by major category\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Assignment # 9 Part 2\n# Replace the single \"College major by Pay.\" paragraph with four new\n# paragraphs of analysis notes, keeping the \"_GoBack\" bookmark anchored\n# right before the final \".\" of the last paragraph.\n\n$d = $word.ActiveDocument\n\n# The document's only (first) paragraph currently holds \"College major by Pay.\"\n$firstRange = $d.Paragraphs.First.Range\n\n# Insert the three new leading paragraphs directly above it, in document order.\n$firstRange.InsertBefore(\"Anova analysis of median pay by Major Category & Unemployment by major category.`rTukey SD based on Anova.`rBayes prediction of mean pay and unemployment by major category.`r\")\n\n# The original paragraph (with \"C\" ... bookmark ... \"ollege major by Pay.\") is\n# now the 4th paragraph. Trim its trailing text (after the bookmark) down to\n# just \".\" and swap its leading text (before the bookmark) for the new\n# \"Linear regression ...\" sentence - preserving the bookmark's position.\n$targetRange = $d.Paragraphs.Item(4).Range\n$find = $targetRange.Find\n$find.Text = \"ollege major by Pay.\"\n$find.MatchCase = $true\n$find.Replacement.Text = \".\"\n$find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n\n$targetRange2 = $d.Paragraphs.Item(4).Range\n$find2 = $targetRange2.Find\n$find2.Text = \"C\"\n$find2.MatchCase = $true\n$find2.Replacement.Text = \"Linear regression of median pay and unemployment by major category\"\n$find2.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n"}
